# Edit slide 1, shape 2 ("Subtitle 2" placeholder) to match the target revision:
#  - reposition / resize the shape
#  - switch body font from Georgia to Aptos
#  - insert a new run "222407776/autunm1455csb32" right after "REGISTER NO AND NMID:"
#  - change "MADRAS UNIVERSITY" -> "university of madras"
#
# Original text frame (paragraphs, 1-based Characters() offsets):
#   1-26   "STUDENT NAME:  KANIMOZHI.H"
#   28-48  "REGISTER NO AND NMID:"
#   50-82  "DEPARTMENT:  BSC COMPUTER SCIENCE"
#   84-153 "COLLEGE/UNIVERSITY:  TAGORE COLLEGE OF ARTS&SCIENCE/ MADRAS UNIVERSITY"
# (positions 27/49/83 are the paragraph-break carriage returns)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

# --- Resize / reposition (EMU -> points, 1 pt = 12700 EMU) ---
$sh.Left   = 1440090 / 12700
$sh.Top    = 2727765 / 12700
$sh.Width  = 8825658 / 12700
$sh.Height = 2890715 / 12700

$tr = $sh.TextFrame.TextRange

# --- Insert the new run into paragraph 2, right after "REGISTER NO AND NMID:" ---
$para2 = $tr.Characters(28, 21)
$null = $para2.InsertAfter("222407776/autunm1455csb32")
$newRun = $tr.Characters(28 + 21, 25)
$newRun.Font.Caps = 0

# --- Update paragraph 4 text (MADRAS UNIVERSITY -> university of madras); the new
#     run above shifted everything after it forward by 25 characters. ---
$para4 = $tr.Characters(84 + 25, 70)
$para4.Text = "COLLEGE/UNIVERSITY:  TAGORE COLLEGE OF ARTS&SCIENCE/ university of madras"

# --- Swap the body font Georgia -> Aptos across the whole text frame ---
$tr.Font.Name = "Aptos"
